# TSK#201: Finish the update-report command
#
# The "survey_data" sheet is renamed to "config" (keeping its physical
# identity / sheetId), a brand-new "survey_data" sheet is inserted right
# after it, and the original survey data is copied verbatim into that new
# sheet. The (now empty) "config" sheet is then populated with the grouping
# configuration that used to live implicitly in code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Snapshot the existing survey_data grid (A1:X6) before we touch
#    anything - read with Value2 (Value is not reliable for reads in this
#    host) so we can re-create it untouched in the new sheet afterwards.
# ---------------------------------------------------------------------
$oldSurvey = $wb.Worksheets.Item("survey_data")

$lastRow = 6
$lastCol = 24  # X

$snapshot = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $oldSurvey.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $key = "$r,$c"
            $snapshot[$key] = $val
        }
    }
}

# ---------------------------------------------------------------------
# 2. Rename the existing sheet to "config" and insert a fresh sheet named
#    "survey_data" right after it.
# ---------------------------------------------------------------------
$oldSurvey.Name = "config"
$config = $wb.Worksheets.Item("config")

$newSurvey = $wb.Worksheets.Add($null, $config)
$newSurvey.Name = "survey_data"

# ---------------------------------------------------------------------
# 3. Re-populate the new "survey_data" sheet with the original data.
# ---------------------------------------------------------------------
foreach ($key in $snapshot.Keys) {
    $parts = $key.Split(",")
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $newSurvey.Cells.Item($r, $c).Value = $snapshot[$key]
}

# ---------------------------------------------------------------------
# 4. Clear out the old data left behind in "config" and write the new
#    configuration grid.
# ---------------------------------------------------------------------
$config.Range("A1:X6").ClearContents()

# Header row (field names used by the update-report command).
$config.Cells.Item(1, 1).Value = "class_name"
$config.Cells.Item(1, 2).Value = "target_group_size"
$config.Cells.Item(1, 3).Value = "target_plus_one_allowed"
$config.Cells.Item(1, 4).Value = "target_minus_one_allowed"
$config.Cells.Item(1, 5).Value = "grouping_passes"
$config.Cells.Item(1, 6).Value = "availability_values_delimiter"
$config.Cells.Item(1, 7).Value = "student_id_field_name"
$config.Cells.Item(1, 8).Value = "timezone_field_name"
$config.Cells.Item(1, 9).Value = "preferred_students_field_names"
$config.Cells.Item(1, 10).Value = "disliked_students_field_names"
$config.Cells.Item(1, 11).Value = "availability_field_names"
$config.Cells.Item(1, 12).Value = "show_preferred_students"
$config.Cells.Item(1, 13).Value = "show_disliked_students"
$config.Cells.Item(1, 14).Value = "show_availability_overlap"
$config.Cells.Item(1, 15).Value = "show_scores"
$config.Cells.Item(1, 16).Value = "output_student_name"
$config.Cells.Item(1, 17).Value = "output_student_email"
$config.Cells.Item(1, 18).Value = "output_student_login"

# Row 2: configured values.
$config.Cells.Item(2, 1).Value = "SER401"
$config.Cells.Item(2, 2).Value = 2
$config.Cells.Item(2, 3).Value = $true
$config.Cells.Item(2, 4).Value = $false
$config.Cells.Item(2, 5).Value = 10
$config.Cells.Item(2, 6).Value = ";,"
$config.Cells.Item(2, 7).Value = "Please select your ASURITE ID"
$config.Cells.Item(2, 8).Value = "In what time zone do you live or will you be during the session? Please use UTC so we can match it easier."
$config.Cells.Item(2, 9).Value = "Preferred team member 1"
$config.Cells.Item(2, 10).Value = "Non-preferred student 1"
$config.Cells.Item(2, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [0:00 AM - 3:00 AM]"
$config.Cells.Item(2, 12).Value = $true
$config.Cells.Item(2, 13).Value = $true
$config.Cells.Item(2, 14).Value = $true
$config.Cells.Item(2, 15).Value = $true
$config.Cells.Item(2, 16).Value = $false
$config.Cells.Item(2, 17).Value = $true
$config.Cells.Item(2, 18).Value = $true

# Additional entries for the multi-value field-name lists (columns I, J, K).
$config.Cells.Item(3, 9).Value = "Preferred team member 2"
$config.Cells.Item(3, 10).Value = "Non-preferred student 2"
$config.Cells.Item(3, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [3:00 AM - 6:00 AM]"

$config.Cells.Item(4, 9).Value = "Preferred team member 3"
$config.Cells.Item(4, 10).Value = "Non-preferred student 3"
$config.Cells.Item(4, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [6:00 AM - 9:00 AM]"

$config.Cells.Item(5, 9).Value = "Preferred team member 4"
$config.Cells.Item(5, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [9:00 AM - 12:00 PM]"

$config.Cells.Item(6, 9).Value = "Preferred team member 5"
$config.Cells.Item(6, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [12:00 PM - 3:00 PM]"

$config.Cells.Item(7, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [3:00 PM - 6:00 PM]"
$config.Cells.Item(8, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [6:00 PM - 9:00 PM]"
$config.Cells.Item(9, 11).Value = "Please choose times that are good for your team to meet. Times are in the Phoenix, AZ time zone! [9:00 PM - 12:00 PM]"
